# Rename the category headers used for the charts (Spanish labels) and add
# a new "Deudas" (Debts) category column (M) with its monthly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel existing headers, B:L, and add new header M1 ---
$ws.Range("A1").Value = "Months"
$ws.Range("B1").Value = "Salario"
$ws.Range("C1").Value = "Comisiones"
$ws.Range("D1").Value = "Ventas"
$ws.Range("E1").Value = "Otros Ingresos"
$ws.Range("F1").Value = "Otros"
$ws.Range("G1").Value = "Domicilio"
$ws.Range("I1").Value = "Entretenimiento"
$ws.Range("J1").Value = "Servicios"
$ws.Range("K1").Value = "Higiene"
$ws.Range("L1").Value = "Seguros"
$ws.Range("M1").Value = "Deudas"
$ws.Range("H1").Value = "Transporte"

# --- New data column M (Deudas) for each month row (2-13) ---
$mValues = @(18000, 30000, 13000, 80000, 24500, 15000, 5000, 5000, 15000, 24500, 5000, 24500)
for ($i = 0; $i -lt $mValues.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 13)
    $cell.Value = $mValues[$i]
    $cell.NumberFormat = $ws.Cells.Item($row, 2).NumberFormat
}

# --- Update the on-screen selection to match the saved view state ---
$ws.Range("H1").Select()
